$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.152.36"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.37%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.862.31"
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'0.9996"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'233.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -3.54%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9996"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.01%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4655"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2822"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -3.13%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06547"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.62%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'20.10"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -0.79%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07805"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.56%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'96.25"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -7.75%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.856.52"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -4.43%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  -3.41%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.6685"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -4.00%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'281.17"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -5.38%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'30.173.96"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -3.21%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.9994"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.08%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'5.467"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.93%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'12.62"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -2.98%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'2.100.36"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -4.87%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.000007245"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -4.88%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.9989"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.03%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'6.140"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -4.73%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'9.301"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -2.77%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'165.14"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -2.25%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  -4.95%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'1.900"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -9.84%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'1.340"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -4.46%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.09597"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -4.75%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'4.389"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -5.32%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'1.467"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -4.55%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'4.107"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -5.66%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.04645"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -3.87%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.6991"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -5.64%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'1.093"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -3.84%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.64%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -5.48%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'6.289"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -9.00%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'2.516"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -5.05%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'72.62"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -5.46%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.8518"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -2.53%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'1.921"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -5.56%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.9992"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.06%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.4154"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -5.21%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'103.27"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.74%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'988.70"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -3.65%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'7.165"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -5.61%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'9.195"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.48%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -3.13%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -5.89%  "
$ws.Range('E51').Style = 'Normal'
